# Applies the "Error Calculations and Plots" data-fix edit:
#  - some previously-missing (blank) cells get their numeric values filled in
#  - a few cells that had values are now re-blanked
#  - the "RM 232" row and the "SC 92" row are removed entirely, shifting all
#    following rows up by two and shrinking the used range to A1:F33

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Per-cell value fixes on rows 2-25 (row numbers unaffected by the later deletes) ---

# C3: 11.2 -> blank
$ws.Range("C3").Value = ""

# D5: -14.4 -> blank
$ws.Range("D5").Value = ""

# E6: blank -> -5.7
$ws.Range("E6").Value = -5.7

# D8: blank -> -13.9
$ws.Range("D8").Value = -13.9

# D10: blank -> -14.7
$ws.Range("D10").Value = -14.7

# E11: blank -> -7.9
$ws.Range("E11").Value = -7.9

# D12: -14.1 -> blank ; E12: -5.3 -> blank
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""

# E13: blank -> -5.3
$ws.Range("E13").Value = -5.3

# D15: blank -> -15.2
$ws.Range("D15").Value = -15.2

# E17: -7.3 -> blank
$ws.Range("E17").Value = ""

# D18: -15.2 -> blank
$ws.Range("D18").Value = ""

# D19: -15.5 -> blank ; E19: -6.5 -> blank
$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""

# D25: blank -> -15.5 ; E25: blank -> -7.1
$ws.Range("D25").Value = -15.5
$ws.Range("E25").Value = -7.1

# --- Remove the "RM 232" row (26) and the "SC 92" row (originally 28, becomes
#     27 once row 26 is removed). Everything below shifts up by one each time. ---

$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# --- Final value fixes on the now-shifted rows (26-33) ---

# Row 26 = "SC 5": blank -> -20.2
$ws.Range("B26").Value = -20.2

# Row 27 = "SC 101": -20.4 -> blank
$ws.Range("B27").Value = ""

# Row 29 = "SC 119": -13 -> blank
$ws.Range("D29").Value = ""

# Row 31 = "SC 132": -8.1 -> blank
$ws.Range("E31").Value = ""

# Row 32 = "SC 193": -6.4 -> blank
$ws.Range("E32").Value = ""

# Row 33 = "SC 232": blank -> -19.5 ; blank -> 10.4
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
